# Update the Cthrc1-Fzd3 LR-pair sheet with new TPM-based values.
#
# Source clusters ECs/FAPs/MuSCs are now combined with only 4 target
# clusters (ECs, FAPs, MuSCs, Neutrophils) instead of 5 - the
# "Resolving-Mac" target-cluster combinations are dropped entirely, which
# also drops the now-unused "Resolving-Mac" shared string on save.
#
# Rows 2-13 (3 source clusters x 4 target clusters) are rewritten in place
# with the new TPM-derived numbers, and the trailing rows 14-16 (which
# held the old "MuSCs x Resolving-Mac" block, no longer present) are
# deleted outright.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One element per column A:T, one row per array entry below.
$newData = @(
    @("ECs","Cthrc1","Fzd3","ECs",3,1,0.63867,1.91601,0.02162993170452444,0.02162993170452444,3,1,0.2310223333333334,0.6930670000000001,0.1371162007804382,0.1371162007804382,0.14754703363,1.32792330267,0.002965814058464739,0.00296581405846474),
    @("ECs","Cthrc1","Fzd3","FAPs",3,1,0.63867,1.91601,0.02162993170452444,0.02162993170452444,3,1,0.4828523333333334,1.448557,0.2865821521640898,0.2865821521640898,0.30838329973,2.77544969757,0.006198752379044893,0.006198752379044895),
    @("ECs","Cthrc1","Fzd3","MuSCs",3,1,0.63867,1.91601,0.02162993170452444,0.02162993170452444,3,1,0.9450603333333335,2.835181,0.5609114952015948,0.5609114952015947,0.60358168309,5.432235147810001,0.01213247733349318,0.01213247733349318),
    @("ECs","Cthrc1","Fzd3","Neutrophils",3,1,0.63867,1.91601,0.02162993170452444,0.02162993170452444,1,0.3333333333333333,0.02593033333333333,0.077791,0.01539015185387714,0.01539015185387714,0.01656092599,0.14904833391,0.0003328879335216227,0.0003328879335216228),
    @("FAPs","Cthrc1","Fzd3","ECs",3,1,28.021538,84.06461399999999,0.9490095874171892,0.9490095874171893,3,1,0.2310223333333334,0.6930670000000001,0.1371162007804382,0.1371162007804382,6.473601092348666,58.262409831138,0.1301245891308562,0.1301245891308562),
    @("FAPs","Cthrc1","Fzd3","FAPs",3,1,28.021538,84.06461399999999,0.9490095874171892,0.9490095874171893,3,1,0.4828523333333334,1.448557,0.2865821521640898,0.2865821521640898,13.53026500688867,121.772385061998,0.271969209986373,0.2719692099863731),
    @("FAPs","Cthrc1","Fzd3","MuSCs",3,1,28.021538,84.06461399999999,0.9490095874171892,0.9490095874171893,3,1,0.9450603333333335,2.835181,0.5609114952015948,0.5609114952015947,26.48204404279267,238.338396385134,0.5323103866388241,0.5323103866388241),
    @("FAPs","Cthrc1","Fzd3","Neutrophils",3,1,28.021538,84.06461399999999,0.9490095874171892,0.9490095874171893,1,0.3333333333333333,0.02593033333333333,0.077791,0.01539015185387714,0.01539015185387714,0.7266078208526665,6.539470387673999,0.01460540166113584,0.01460540166113584),
    @("MuSCs","Cthrc1","Fzd3","ECs",3,1,0.866931,2.600793,0.02936048087828625,0.02936048087828625,3,1,0.2310223333333334,0.6930670000000001,0.1371162007804382,0.1371162007804382,0.200280422459,1.802523802131,0.004025797591117314,0.004025797591117314),
    @("MuSCs","Cthrc1","Fzd3","FAPs",3,1,0.866931,2.600793,0.02936048087828625,0.02936048087828625,3,1,0.4828523333333334,1.448557,0.2865821521640898,0.2865821521640898,0.418599656189,3.767396905701,0.008414189798671879,0.008414189798671879),
    @("MuSCs","Cthrc1","Fzd3","MuSCs",3,1,0.866931,2.600793,0.02936048087828625,0.02936048087828625,3,1,0.9450603333333335,2.835181,0.5609114952015948,0.5609114952015947,0.8193020998370002,7.373718898533001,0.01646863122927737,0.01646863122927737),
    @("MuSCs","Cthrc1","Fzd3","Neutrophils",3,1,0.866931,2.600793,0.02936048087828625,0.02936048087828625,1,0.3333333333333333,0.02593033333333333,0.077791,0.01539015185387714,0.01539015185387714,0.022479809807,0.202318288263,0.0004518622592196815,0.0004518622592196815)
)

$startRow = 2
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $rowValues = $newData[$i]
    for ($c = 0; $c -lt $rowValues.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value2 = $rowValues[$c]
    }
}

# Drop the old trailing rows that belonged to the now-removed
# "Resolving-Mac" target cluster.
$ws.Rows("14:16").Delete()
